$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.508.73"
$ws.Range("E2").Value = "  +2.11%  "
$ws.Range("D3").Value = "1.869.00"
$ws.Range("E3").Value = "  +1.21%  "
$ws.Range("E4").Value = "  +0.34%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "312.32"
$ws.Range("D5").Style = "Normal"
$ws.Range("E6").Value = "  +0.38%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4781"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.19%  "
$ws.Range("E8").Value = "  +2.89%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07353"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +1.94%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.9372"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +1.04%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "20.71"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +5.02%  "
$ws.Range("E12").Value = "  +2.17%  "
$ws.Range("D13").Value = "1.831.06"
$ws.Range("E13").Value = "  -2.34%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.439"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +2.28%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.573"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +2.61%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "90.75"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +2.09%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.014"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.30%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000008892"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +2.96%  "
$ws.Range("E19").Value = "  +0.34%  "
$ws.Range("E20").Value = "  +2.59%  "
$ws.Range("D21").Value = "27.514.01"
$ws.Range("E21").Value = "  +2.01%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.130"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.51%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "10.73"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.86%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.958"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.90%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "153.87"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.80%  "
$ws.Range("E26").Value = "  +1.89%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.028"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +1.59%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "115.88"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.39%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "5.000"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +1.59%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.08924"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.54%  "
$ws.Range("E31").Value = "  +0.80%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.216"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +3.44%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.610"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +2.69%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.7530"
$ws.Range("D34").Style = "Normal"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.712"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.39%  "
$ws.Range("E36").Value = "  +4.84%  "
$ws.Range("E37").Value = "  +1.14%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.005"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.99%  "
$ws.Range("E39").Value = "  +0.23%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.5348"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +2.93%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "7.085"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +1.91%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1526"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +1.10%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "8.487"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +3.55%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "10.70"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +1.29%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.4810"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +1.86%  "
$ws.Range("E46").Value = "  +0.33%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.661"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +3.64%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "102.85"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.35%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "67.48"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +2.18%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.06086"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +1.12%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.9240"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +4.39%  "
